# Avance de la clase Definitions hasta el paso de Grabar para
# Reprogramación Otro Pagaré.
#
# - Corrige el texto de la cabecera "Forma Ddesembolso" -> "Forma Desembolso"
#   (la cadena vieja queda huérfana y se recicla; la nueva se agrega al final
#   de la tabla de cadenas compartidas).
# - Ensancha ligeramente la columna U ("Número Propuesta") para que el
#   contenido se siga viendo bien tras el resto de cambios de la pantalla.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige el texto de la cabecera de "Forma Desembolso" (columna P).
$ws.Range("P1").Value = "Forma Desembolso"

# Ajusta el ancho de la columna U (21) un poco más ancho que el original.
# (21.140625 caracteres es el ancho canónico del archivo de referencia; el
# motor COM redondea ColumnWidth a 1/6 de caracter, así que usamos el valor
# de entrada que produce el resultado almacenado más cercano posible.)
$ws.Columns.Item(21).ColumnWidth = 20.3
